# Cap nhat xuat bao cao voi kha nang tu format
#
# The template's title cell (A1) used to hold the literal text "A1"
# (a left-over placeholder). Clear it so the cell is blank, and leave the
# active selection on B1, matching the saved worksheet state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "A1" placeholder text from cell A1 (keep its formatting).
$ws.Range("A1").ClearContents()

# Record B1 as the selected/active cell, as it was when the sheet was saved.
$ws.Range("B1").Select()
